# Tyler_Kroft_TE_2018.xlsx - scraper update
# Adds two new columns ("height" and "weight") between the existing
# "fumbles" column (D) and "fantasy points" column (E), pushing the
# "fantasy points" data from column E to column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "fantasy points" column (E) two columns to the right
# (to G), opening up blank columns E and F for the new fields.
$ws.Range("E1:E6").Insert(-4161)   # xlShiftToRight
$ws.Range("F1:F6").Insert(-4161)   # xlShiftToRight

# New column headers
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# New column data (same value repeated for every player row, as scraped)
$ws.Range("E2:E6").Value = 6.5
$ws.Range("F2:F6").Value = 252
